# Fix the "Scuidal " typo to "Suicidal" in the category column (column A)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Suicidal"
$ws.Range("A12").Value = "Suicidal"
